$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Company UEN" column values (D2:D4) from "34242" to "321432".
# Leading apostrophe preserves the existing text/quote-prefix cell style.
$ws.Range("D2").Value = "'321432"
$ws.Range("D3").Value = "'321432"
$ws.Range("D4").Value = "'321432"

# Update the active selection to match the saved workbook state
$ws.Range("F9:F10").Select()
